# Update "想去人数" (interested count) figures for several events.
# Sheet "展览" (exhibitions) and the consolidated "全部类型" (all types)
# sheet both list the same events, so both need the matching updates.

$wb = $excel.ActiveWorkbook

$wsExhibitions = $wb.Worksheets.Item("展览")
$wsExhibitions.Range("F4").Value = 50
$wsExhibitions.Range("F5").Value = 4986
$wsExhibitions.Range("F7").Value = 84
$wsExhibitions.Range("F8").Value = 292
$wsExhibitions.Range("F9").Value = 46

$wsAllTypes = $wb.Worksheets.Item("全部类型")
$wsAllTypes.Range("F8").Value = 50
$wsAllTypes.Range("F9").Value = 4986
$wsAllTypes.Range("F11").Value = 84
$wsAllTypes.Range("F13").Value = 292
$wsAllTypes.Range("F14").Value = 46
